$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 32,2
$arr[0,0] = -0.098389817566534532
$arr[0,1] = 0.097855647603999785
$arr[1,0] = -0.11675307120240586
$arr[1,1] = 0.11519165509362672
$arr[2,0] = -0.0654747902382109
$arr[2,1] = 0.065009460279076947
$arr[3,0] = -0.057009460361811648
$arr[3,1] = 0.056591213643544691
$arr[4,0] = -0.053591213686492445
$arr[4,1] = 0.052172202290058145
$arr[5,0] = -0.029879679249052771
$arr[5,1] = 0.029674544341085962
$arr[6,0] = -0.019674544452588538
$arr[6,1] = 0.019638709559127943
$arr[7,0] = -0.0096387096727754695
$arr[7,1] = 0.0096064900667305331
$arr[8,0] = -0.0076064901128365392
$arr[8,1] = 0.0075883823465239253
$arr[9,0] = -0.0055883823937055155
$arr[9,1] = 0.0055888680594122775
$arr[10,0] = -0.0025888681152430593
$arr[10,1] = 0.0025890515741897602
$arr[11,0] = 0.00091094836564531434
$arr[11,1] = -0.00091131563049362541
$arr[12,0] = 0.0044113155706009266
$arr[12,1] = -0.0044129546328859348
$arr[13,0] = -0.0090809763536645605
$arr[13,1] = 0.0090526289093704193
$arr[14,0] = -0.0080526289480937763
$arr[14,1] = 0.0080343481509448367
$arr[15,0] = -0.0060343481991269599
$arr[15,1] = 0.0060032422899398696
$arr[16,0] = -0.0040032423392464267
$arr[16,1] = 0.0039999999331445935
$arr[17,0] = 0.020752410149125211
$arr[17,1] = -0.020859615454050129
$arr[18,0] = -0.012091173233750752
$arr[18,1] = 0.012016121845738414
$arr[19,0] = -0.0080161218838092907
$arr[19,1] = 0.0080056098545249199
$arr[20,0] = -0.0040056098930270068
$arr[20,1] = 0.0039999999612545523
$arr[21,0] = -0.045716864900272824
$arr[21,1] = 0.0455023819647149
$arr[22,0] = -0.040502382021103678
$arr[22,1] = 0.040099595994467663
$arr[23,0] = -0.020099596185700896
$arr[23,1] = 0.019999999806175062
$arr[24,0] = -0.020292523082805047
$arr[24,1] = 0.020229851945405386
$arr[25,0] = -0.017729851992063672
$arr[25,1] = 0.01765164289040122
$arr[26,0] = -0.015151642937763832
$arr[26,1] = 0.014703431206426032
$arr[27,0] = -0.01761409557868987
$arr[27,1] = 0.017312812457054072
$arr[28,0] = -0.010312812548319172
$arr[28,1] = 0.01023641759395133
$arr[29,0] = -0.021170866717706183
$arr[29,1] = 0.021022378321646329
$arr[30,0] = -0.014022378417227088
$arr[30,1] = 0.014000941059537197
$arr[31,0] = -0.0040009411816139817
$arr[31,1] = 0.0039999999295616817

$ws.Range("A1:B32").Value = $arr

# Column widths (closest achievable via ColumnWidth, which this runtime
# quantizes to 1/6-character steps using a fixed MDW=6 model):
# target stored width A = 15.7109375  -> nearest achievable 15.666666666666666
# target stored width B = 16.42578125 -> nearest achievable 16.5
$ws.Columns.Item(1).ColumnWidth = 14.833333333333334
$ws.Columns.Item(2).ColumnWidth = 15.666666666666666

